$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws.Range("D2") '305.21'
Set-TextValue $ws.Range("E2") '0.50%'
Set-TextValue $ws.Range("D3") '36.17'
Set-TextValue $ws.Range("E3") '-3.04%'
Set-TextValue $ws.Range("D4") '5.022'
Set-TextValue $ws.Range("E4") '-0.31%'
Set-TextValue $ws.Range("D5") '0.07879'
Set-TextValue $ws.Range("E5") '0.42%'
Set-TextValue $ws.Range("D6") '2.129'
Set-TextValue $ws.Range("E6") '-3.65%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range("D7") '7.947'
Set-TextValue $ws.Range("E7") '-0.48%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D8") '4.130'
Set-TextValue $ws.Range("E8") '2.58%'
Set-TextValue $ws.Range("D9") '0.9228'
Set-TextValue $ws.Range("E9") '-0.35%'
Set-TextValue $ws.Range("D10") '0.09715'
Set-TextValue $ws.Range("E10") '-1.50%'
Set-TextValue $ws.Range("D11") '0.1847'
Set-TextValue $ws.Range("E11") '-1.86%'
Set-TextValue $ws.Range("D12") '0.08655'
Set-TextValue $ws.Range("E12") '-0.67%'
Set-TextValue $ws.Range("D13") '0.03559'
Set-TextValue $ws.Range("E13") '-1.51%'
Set-TextValue $ws.Range("D14") '0.09911'
Set-TextValue $ws.Range("E14") '-0.35%'
Set-TextValue $ws.Range("D15") '0.001440'
Set-TextValue $ws.Range("E15") '-2.78%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range("D16") '0.04584'
Set-TextValue $ws.Range("E16") '-0.38%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D17") '0.005698'
Set-TextValue $ws.Range("E17") '0.88%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D18") '3.471'
Set-TextValue $ws.Range("E18") '0.29%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D19") '2.753'
Set-TextValue $ws.Range("E19") '17.33%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range("D20") '0.3375'
Set-TextValue $ws.Range("E20") '-1.71%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D21") '0.1349'
Set-TextValue $ws.Range("E21") '1.68%'
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range("D22") '5.155'
Set-TextValue $ws.Range("E22") '7.99%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range("D23") '0.2208'
Set-TextValue $ws.Range("E23") '0.27%'
Set-TextValue $ws.Range("D24") '0.001232'
Set-TextValue $ws.Range("E24") '-1.60%'
Set-TextValue $ws.Range("D25") '0.004806'
Set-TextValue $ws.Range("E25") '-7.66%'
Set-TextValue $ws.Range("D26") '0.0001300'
Set-TextValue $ws.Range("E26") '-7.20%'
Set-TextValue $ws.Range("D27") '0.0004750'
Set-TextValue $ws.Range("E27") '74.69%'
Set-TextValue $ws.Range("D39") '0.01854'
Set-TextValue $ws.Range("E39") '1.34%'
Set-TextValue $ws.Range("D40") '0.04711'
Set-TextValue $ws.Range("E40") '-0.92%'
Set-TextValue $ws.Range("D41") '0.007796'
Set-TextValue $ws.Range("E41") '-1.81%'
Set-TextValue $ws.Range("D42") '0.1383'
Set-TextValue $ws.Range("E42") '-1.36%'
Set-TextValue $ws.Range("D43") '0.007716'
Set-TextValue $ws.Range("E43") '2.01%'
Set-TextValue $ws.Range("E44") '-3.54%'
Set-TextValue $ws.Range("D45") '0.01137'
Set-TextValue $ws.Range("E45") '9.48%'
Set-TextValue $ws.Range("D46") '0.00006376'
Set-TextValue $ws.Range("E46") '1.10%'
Set-TextValue $ws.Range("D47") '0.00000000750'
Set-TextValue $ws.Range("E47") '-0.04%'
Set-TextValue $ws.Range("E48") '0.13%'
Set-TextValue $ws.Range("D49") '51.93'
Set-TextValue $ws.Range("E49") '44.64%'
Set-TextValue $ws.Range("D50") '0.001900'
Set-TextValue $ws.Range("E50") '-29.38%'
Set-TextValue $ws.Range("D51") '0.00002101'
Set-TextValue $ws.Range("E51") '-0.04%'
